$p = $ppt.ActivePresentation

# --- Slide Master 1: normalize the "outline text" placeholder font sizes ---
# All seven outline levels in the body placeholder are unified to 18pt.
$masterBody = $p.SlideMaster.Shapes.Item(2)
$bodyRange = $masterBody.TextFrame.TextRange
$paraCount = $bodyRange.Paragraphs().Count
for ($i = 1; $i -le $paraCount; $i++) {
    $para = $bodyRange.Paragraphs($i, 1)
    $para.Font.Size = 18
}

# --- Slides 6 & 7: shorten the shared title to "Creating a relationship" ---
$p.Slides.Item(6).Shapes.Item(1).TextFrame.TextRange.Text = "Creating a relationship"
$p.Slides.Item(7).Shapes.Item(1).TextFrame.TextRange.Text = "Creating a relationship"

# --- Slide 8: retitle to "Population" ---
$p.Slides.Item(8).Shapes.Item(1).TextFrame.TextRange.Text = "Population"
